# Generate Report for Handback
# Updates the handback-status report timestamps for the
# "6ea5906b-87f6-4695-b6e9-88e0df925214" file across all sheets, as a fresh
# handback report generation would refresh the handoff/handback timestamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: refresh "Latest HO Xliff Generate Date" for row 3 ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-10-19 11:09:34"

# --- zh-cn sheet: refresh Correspond Handoff/Handback datetimes for row 3 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-10-19 11:09:22"
$zhcn.Range("K3").Value = "2016-10-19 11:10:06"

# --- de-de sheet: refresh Correspond Handoff/Handback datetimes for row 3 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-10-19 11:09:34"
$dede.Range("K3").Value = "2016-10-19 11:10:24"
